$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where E column should become a numeric 100 (was inline string "100")
$numericRows = @(2, 3, 4, 5, 6, 7, 8, 11, 12)
foreach ($r in $numericRows) {
    $ws.Cells.Item($r, 5).Value = 100
}

# Row 9 & 10: clear E (ref2/ref1) and Q/R/S (Common_name/Species_name/Category)
$ws.Cells.Item(9, 5).ClearContents()
$ws.Cells.Item(9, 17).ClearContents()
$ws.Cells.Item(9, 18).ClearContents()
$ws.Cells.Item(9, 19).ClearContents()

$ws.Cells.Item(10, 5).ClearContents()
$ws.Cells.Item(10, 17).ClearContents()
$ws.Cells.Item(10, 18).ClearContents()
$ws.Cells.Item(10, 19).ClearContents()
